$wb = $excel.ActiveWorkbook

$oldGuid = "1c3f920e-ebcc-4b8b-886e-c6571d934d3d"
$newGuid = "f0d15859-f377-4e16-94c0-329d2d4de750"
$oldHash = "129b96a823794634001c1ab2703ec28f506745dc"
$newHash = "5a7cd9757f250c791862152d2b7ed58834670876"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: file name, path and name (hyperlink display text), generate date
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-22 12:59:54"

# zh-cn sheet: file name (hyperlink display text), handoff file, handoff datetime
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-22 12:59:49"

# de-de sheet: file name (hyperlink display text), handoff file, handoff datetime
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-22 12:59:54"
